# This script applies a cyclic rotation of the species-observation data
# held in rows 2, 3 and 5 of the active worksheet:
#   row 2 <- (old) row 5
#   row 3 <- (old) row 2
#   row 5 <- (old) row 3
# Columns C, H, I, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY
# are identical across these three rows and are left untouched.
#
# Note: use .Value2 (not .Value) to read raw/unformatted values so numbers
# and strings round-trip correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the "before" values for the columns that differ between rows ---

$row2A  = $ws.Range("A2").Value2
$row2B  = $ws.Range("B2").Value2
$row2D  = $ws.Range("D2").Value2
$row2E  = $ws.Range("E2").Value2
$row2F  = $ws.Range("F2").Value2
$row2G  = $ws.Range("G2").Value2
$row2P  = $ws.Range("P2").Value2
$row2Q  = $ws.Range("Q2").Value2
$row2R  = $ws.Range("R2").Value2
$row2AI = $ws.Range("AI2").Value2

$row3A  = $ws.Range("A3").Value2
$row3B  = $ws.Range("B3").Value2
$row3D  = $ws.Range("D3").Value2
$row3E  = $ws.Range("E3").Value2
$row3F  = $ws.Range("F3").Value2
$row3G  = $ws.Range("G3").Value2
$row3P  = $ws.Range("P3").Value2
$row3Q  = $ws.Range("Q3").Value2
$row3R  = $ws.Range("R3").Value2
$row3AI = $ws.Range("AI3").Value2

$row5A  = $ws.Range("A5").Value2
$row5B  = $ws.Range("B5").Value2
$row5D  = $ws.Range("D5").Value2
$row5E  = $ws.Range("E5").Value2
$row5F  = $ws.Range("F5").Value2
$row5G  = $ws.Range("G5").Value2
$row5P  = $ws.Range("P5").Value2
$row5Q  = $ws.Range("Q5").Value2
$row5R  = $ws.Range("R5").Value2
$row5AI = $ws.Range("AI5").Value2

# --- Write row 2 <- old row 5 ---
$ws.Range("A2").Value = $row5A
$ws.Range("B2").Value = $row5B
$ws.Range("D2").Value = $row5D
$ws.Range("E2").Value = $row5E
$ws.Range("F2").Value = $row5F
$ws.Range("G2").Value = $row5G
$ws.Range("P2").Value = $row5P
$ws.Range("Q2").Value = $row5Q
$ws.Range("R2").Value = $row5R
$ws.Range("AI2").ClearContents()

# --- Write row 3 <- old row 2 ---
$ws.Range("A3").Value = $row2A
$ws.Range("B3").Value = $row2B
$ws.Range("D3").Value = $row2D
$ws.Range("E3").Value = $row2E
$ws.Range("F3").Value = $row2F
$ws.Range("G3").Value = $row2G
$ws.Range("P3").Value = $row2P
$ws.Range("Q3").Value = $row2Q
$ws.Range("R3").Value = $row2R
$ws.Range("AI3").Value = $row2AI

# --- Write row 5 <- old row 3 ---
$ws.Range("A5").Value = $row3A
$ws.Range("B5").Value = $row3B
$ws.Range("D5").Value = $row3D
$ws.Range("E5").Value = $row3E
$ws.Range("F5").Value = $row3F
$ws.Range("G5").Value = $row3G
$ws.Range("P5").Value = $row3P
$ws.Range("Q5").Value = $row3Q
$ws.Range("R5").Value = $row3R
$ws.Range("AI5").Value = $row3AI
